$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 4

$ws.Cells.Item($row, 1).Value = "Tue Oct 31 00:50:03 2023"
$ws.Cells.Item($row, 2).Value = "SZ399006"
$ws.Cells.Item($row, 3).Value = "创业板指"
$ws.Cells.Item($row, 4).Value = "未开盘"
$ws.Cells.Item($row, 5).NumberFormat = "@"
$ws.Cells.Item($row, 5).Value = "1977.69"
$ws.Cells.Item($row, 6).Value = "+47.78  +2.48%"
$ws.Cells.Item($row, 7).Value = 1983.46
$ws.Cells.Item($row, 8).Value = 1933.79
$ws.Cells.Item($row, 9).Value = 169000000
$ws.Cells.Item($row, 10).Value = 0.036
$ws.Cells.Item($row, 11).Value = 1931.7
$ws.Cells.Item($row, 12).Value = 1929.91
$ws.Cells.Item($row, 13).Value = 277975000000
$ws.Cells.Item($row, 14).Value = 1.14
$ws.Cells.Item($row, 15).Value = 2661.28
$ws.Cells.Item($row, 16).Value = 765
$ws.Cells.Item($row, 17).Value = 0.027
$ws.Cells.Item($row, 18).Value = 23
$ws.Cells.Item($row, 19).Value = 1840.26
$ws.Cells.Item($row, 20).Value = 163
$ws.Cells.Item($row, 21).Value = 115000
